# Update faturamento_diario_lojas.xlsx (Dados_BIBI_PF) with new AB-column
# values and recalculated AG (total) column for rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New AB values (previously 0)
$ws.Range("AB2").Value = 5032.31
$ws.Range("AB3").Value = 13806
$ws.Range("AB4").Value = 3748
$ws.Range("AB5").Value = 1570
$ws.Range("AB6").Value = 24156.31

# Updated AG totals (sum of B:AF) reflecting the new AB values
$ws.Range("AG2").Value = 295780.85
$ws.Range("AG3").Value = 186925.3
$ws.Range("AG4").Value = 91718.33
$ws.Range("AG5").Value = 81397.48
$ws.Range("AG6").Value = 655821.96
